# Edit the "Hoja3" worksheet per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")

# --- Cell value changes -----------------------------------------------
# M1: 2011 -> 2013
$ws.Range("M1").Value = 2013

# AR1: 256850 -> 234.34
$ws.Range("AR1").Value = 234.34

# AR2: 125820 -> 3456.56
$ws.Range("AR2").Value = 3456.56

# --- Shared string change ----------------------------------------------
# Column P (P1:P145) all shared the single string "RACG770826MBCMSS".
# Re-write the whole block at once so the engine rewrites the single
# shared-string table entry in place instead of allocating a new one.
$ws.Range("P1:P145").Value = "PECJ971112HMCRHN"

# --- Extend P column formatting (style only, no value) to rows 148:152,
#     matching the blank, wrap-text styled cells already present on
#     P146:P147.
$ws.Range("P1").Copy() | Out-Null
$ws.Range("P148:P152").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Sheet view / selection --------------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 108
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("P1:P145").Select() | Out-Null
